$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels in B1 and C1 (SD [mg/l] <-> Conc [mg/l])
$headerB = $ws.Range("B1").Value()
$headerC = $ws.Range("C1").Value()
$ws.Range("B1").Value = $headerC
$ws.Range("C1").Value = $headerB

# Swap the data values between columns B and C for rows 2 through 13
for ($row = 2; $row -le 13; $row++) {
    $valB = $ws.Cells.Item($row, 2).Value()
    $valC = $ws.Cells.Item($row, 3).Value()
    $ws.Cells.Item($row, 2).Value = $valC
    $ws.Cells.Item($row, 3).Value = $valB
}
